$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# row 12
$ws.Range("H12").Value = 271.22223
$ws.Range("I12").Value = 234.57143
$ws.Range("J12").Value = 399.5
$ws.Range("K12").Value = 234.57143
$ws.Range("L12").Value = 399.5
$ws.Range("M12").Value = -64.57142999999999
$ws.Range("N12").Value = -739.5
# row 18
$ws.Range("H18").Value = 780.8182
$ws.Range("I18").Value = 738.9
$ws.Range("J18").Value = 1200
$ws.Range("K18").Value = 738.9
$ws.Range("L18").Value = 1200
$ws.Range("M18").Value = -454.9
$ws.Range("N18").Value = -1768
# row 39
$ws.Range("H39").Value = 1428
$ws.Range("I39").Value = 837.5833
$ws.Range("J39").Value = 3789.6667
$ws.Range("K39").Value = 2512.7499
$ws.Range("L39").Value = 11369.0001
$ws.Range("M39").Value = -2216.7499
$ws.Range("N39").Value = -11961.0001
# row 103
$ws.Range("H103").Value = 2450.2
$ws.Range("I103").Value = 5277
$ws.Range("J103").Value = 1743.5
$ws.Range("K103").Value = 15831
$ws.Range("L103").Value = 5230.5
$ws.Range("M103").Value = -15245
$ws.Range("N103").Value = -6402.5
# row 138
$ws.Range("H138").Value = 2795.2415
# row 141
$ws.Range("H141").Value = 2308.4443
$ws.Range("I141").Value = 2109.625
$ws.Range("J141").Value = 3899
$ws.Range("K141").Value = 6328.875
$ws.Range("L141").Value = 11697
$ws.Range("M141").Value = -1148.875
$ws.Range("N141").Value = -22057

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 32477.383
$ws.Range("I32").Value = 33604.227
$ws.Range("K32").Value = 33604.227
$ws.Range("M32").Value = -33317.227
# row 45
$ws.Range("H45").Value = 3907.923
$ws.Range("I45").Value = 2836.6
$ws.Range("J45").Value = 4577.5
$ws.Range("K45").Value = 2836.6
$ws.Range("L45").Value = 4577.5
$ws.Range("M45").Value = -2459.6
$ws.Range("N45").Value = -5331.5
# row 74
$ws.Range("H74").Value = 3718.88
$ws.Range("I74").Value = 3727.3914
$ws.Range("K74").Value = 3727.3914
$ws.Range("M74").Value = -2853.3914
# row 77
$ws.Range("H77").Value = 3718.88
$ws.Range("I77").Value = 3727.3914
$ws.Range("K77").Value = 18636.957
$ws.Range("M77").Value = -14268.957
# row 97
$ws.Range("H97").Value = 1504.5
$ws.Range("I97").Value = 1130.875
$ws.Range("K97").Value = 1130.875
$ws.Range("M97").Value = -634.875
# row 127
$ws.Range("H127").Value = 35674.31
$ws.Range("I127").Value = 25797.092
$ws.Range("K127").Value = 25797.092
$ws.Range("M127").Value = -20837.092

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
# row 99
$ws.Range("H99").Value = 5175.7
$ws.Range("I99").Value = 4644.8335
$ws.Range("K99").Value = 4644.8335
$ws.Range("M99").Value = -3146.8335

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
# row 99
$ws.Range("H99").Value = 3771.818
$ws.Range("I99").Value = 3913.4285
$ws.Range("K99").Value = 3913.4285
$ws.Range("M99").Value = -2415.4285
# row 126
$ws.Range("H126").Value = 3771.818
$ws.Range("I126").Value = 3913.4285
$ws.Range("K126").Value = 11740.2855
$ws.Range("M126").Value = -9270.2855
# row 132
$ws.Range("H132").Value = 2000
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 6000
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -11060

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
# row 2
$ws.Range("H2").Value = 145.60869
$ws.Range("I2").Value = 100
$ws.Range("J2").Value = 249.85715
$ws.Range("K2").Value = 100
$ws.Range("L2").Value = 249.85715
$ws.Range("M2").Value = 13
$ws.Range("N2").Value = -475.85715
# row 70
$ws.Range("H70").Value = 2627.1333
$ws.Range("I70").Value = 1888.5
$ws.Range("J70").Value = 4104.4
$ws.Range("K70").Value = 1888.5
$ws.Range("L70").Value = 4104.4
$ws.Range("M70").Value = -1618.5
$ws.Range("N70").Value = -4644.4
# row 73
$ws.Range("H73").Value = 2627.1333
$ws.Range("I73").Value = 1888.5
$ws.Range("J73").Value = 4104.4
$ws.Range("K73").Value = 1888.5
$ws.Range("L73").Value = 4104.4
$ws.Range("M73").Value = -952.5
$ws.Range("N73").Value = -5976.4
# row 97
$ws.Range("H97").Value = 1219.4445
$ws.Range("I97").Value = 1026.1538
$ws.Range("J97").Value = 1722
$ws.Range("K97").Value = 1026.1538
$ws.Range("L97").Value = 1722
$ws.Range("M97").Value = -530.1538
$ws.Range("N97").Value = -2714
# row 107
$ws.Range("H107").Value = 63096.5
$ws.Range("I107").Value = 91085.91
$ws.Range("K107").Value = 91085.91
$ws.Range("M107").Value = -89165.91
# row 122
$ws.Range("H122").Value = 3116.5264
$ws.Range("I122").Value = 2485.077
$ws.Range("J122").Value = 4484.6665
$ws.Range("K122").Value = 7455.231000000001
$ws.Range("L122").Value = 13453.9995
$ws.Range("M122").Value = -5005.231000000001
$ws.Range("N122").Value = -18353.9995
# row 132
$ws.Range("H132").Value = 501499.5
$ws.Range("I132").Value = 501499.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 1504498.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1501968.5
$ws.Range("N132").ClearContents()

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
# row 16
$ws.Range("H16").Value = 2083.1538
$ws.Range("I16").Value = 2189.2727
$ws.Range("K16").Value = 2189.2727
$ws.Range("M16").Value = -2019.2727
# row 40
$ws.Range("H40").Value = 5864.913
$ws.Range("I40").Value = 5544.2
$ws.Range("K40").Value = 5544.2
$ws.Range("M40").Value = -5408.2
# row 93
$ws.Range("H93").Value = 2777.5334
$ws.Range("I93").Value = 2773.875
$ws.Range("J93").Value = 2781.7144
$ws.Range("K93").Value = 2773.875
$ws.Range("L93").Value = 2781.7144
$ws.Range("M93").Value = -1525.875
$ws.Range("N93").Value = -5277.7144
# row 100
$ws.Range("H100").Value = 2654.6
$ws.Range("J100").Value = 2548.3333
$ws.Range("L100").Value = 2548.3333
$ws.Range("N100").Value = -3630.3333
# row 122
$ws.Range("H122").Value = 4434.1665
$ws.Range("I122").Value = 4222.5835
$ws.Range("J122").Value = 4857.3335
$ws.Range("K122").Value = 12667.7505
$ws.Range("L122").Value = 14572.0005
$ws.Range("M122").Value = -10217.7505
$ws.Range("N122").Value = -19472.0005
# row 132
$ws.Range("H132").Value = 58501.637
$ws.Range("I132").Value = 69541.05499999999
$ws.Range("K132").Value = 208623.165
$ws.Range("M132").Value = -206093.165
# row 133
$ws.Range("H133").Value = 88343.836
$ws.Range("J133").Value = 88343.836
$ws.Range("L133").Value = 88343.836
$ws.Range("N133").Value = -93403.836

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
# row 62
$ws.Range("H62").Value = 999999
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 999999
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 999999
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -1001247
# row 65
$ws.Range("H65").Value = 999999
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 999999
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 4999995
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -5006235
# row 81
$ws.Range("H81").Value = 3366.8
$ws.Range("I81").Value = 1089
$ws.Range("J81").Value = 3936.25
$ws.Range("K81").Value = 2178
$ws.Range("L81").Value = 7872.5
$ws.Range("M81").Value = -1117
$ws.Range("N81").Value = -9994.5
# row 84
$ws.Range("H84").Value = 3366.8
$ws.Range("I84").Value = 1089
$ws.Range("J84").Value = 3936.25
$ws.Range("K84").Value = 10890
$ws.Range("L84").Value = 39362.5
$ws.Range("M84").Value = -5586
$ws.Range("N84").Value = -49970.5
# row 122
$ws.Range("H122").Value = 2369.2856
$ws.Range("I122").Value = 2197.5833
$ws.Range("K122").Value = 6592.749899999999
$ws.Range("M122").Value = -4142.749899999999
# row 126
$ws.Range("H126").Value = 4400.5835
$ws.Range("I126").Value = 3874.625
$ws.Range("K126").Value = 11623.875
$ws.Range("M126").Value = -9153.875
